$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.429.98'
$ws.Range('E2').Value = '  -0.38%  '

$ws.Range('D3').Value = '3.160.46'
$ws.Range('E3').Value = '  -0.34%  '

$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').Value = '''613.35'
$ws.Range('E5').Value = '  +2.45%  '

$ws.Range('D6').Value = '''148.29'
$ws.Range('E6').Value = '  -2.12%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('D8').Value = '3.156.68'
$ws.Range('E8').Value = '  -0.50%  '

$ws.Range('D9').Value = '''0.527'
$ws.Range('E9').Value = '  -0.62%  '

$ws.Range('E10').Value = '  -0.50%  '

$ws.Range('D11').Value = '''5.46'
$ws.Range('E11').Value = '  -1.66%  '

$ws.Range('E12').Value = '  -0.25%  '

$ws.Range('D13').Value = '''0.0000260'
$ws.Range('E13').Value = '  +0.39%  '

$ws.Range('D14').Value = '''35.62'
$ws.Range('E14').Value = '  -3.32%  '

$ws.Range('D15').Value = '3.668.37'
$ws.Range('E15').Value = '  -0.71%  '

$ws.Range('E16').Value = '  +2.92%  '

$ws.Range('D17').Value = '64.328.57'
$ws.Range('E17').Value = '  -0.58%  '

$ws.Range('D18').Value = '3.158.37'
$ws.Range('E18').Value = '  -0.86%  '

$ws.Range('D19').Value = '''6.92'
$ws.Range('E19').Value = '  -1.30%  '

$ws.Range('D20').Value = '''479.14'
$ws.Range('E20').Value = '  +0.00%  '

$ws.Range('D21').Value = '''14.71'
$ws.Range('E21').Value = '  -0.49%  '

$ws.Range('D22').Value = '''0.715'
$ws.Range('E22').Value = '  +0.37%  '

$ws.Range('D23').Value = '''8.05'
$ws.Range('E23').Value = '  +4.25%  '

$ws.Range('D24').Value = '''13.73'
$ws.Range('E24').Value = '  -0.51%  '

$ws.Range('D25').Value = '''83.95'
$ws.Range('E25').Value = '  -0.06%  '

$ws.Range('E26').Value = '  +0.07%  '

$ws.Range('E27').Value = '  -2.38%  '

$ws.Range('D28').Value = '''8.59'
$ws.Range('E28').Value = '  -0.61%  '

$ws.Range('D29').Value = '''7.12'
$ws.Range('E29').Value = '  +0.93%  '

$ws.Range('D30').Value = '''0.119'
$ws.Range('E30').Value = '  +0.08%  '

$ws.Range('D31').Value = '''2.11'
$ws.Range('E31').Value = '  -6.38%  '

$ws.Range('E32').Value = '  +0.10%  '

$ws.Range('D33').Value = '''2.72'
$ws.Range('E33').Value = '  +0.06%  '

$ws.Range('D34').Value = '''26.38'
$ws.Range('E34').Value = '  -1.62%  '

$ws.Range('E35').Value = '  +2.18%  '

$ws.Range('D36').Value = '0.0₃0803'
$ws.Range('E36').Value = '  +8.87%  '

$ws.Range('D37').Value = '''6.02'
$ws.Range('E37').Value = '  -1.33%  '

$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').Value = '''3.22'
$ws.Range('E38').Value = '  +0.68%  '

$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').Value = '''53.21'
$ws.Range('E39').Value = '  -2.58%  '

$ws.Range('D40').Value = '''463.87'

$ws.Range('D41').Value = '''0.0401'
$ws.Range('E41').Value = '  -0.18%  '

$ws.Range('E42').Value = '  -3.76%  '

$ws.Range('D43').Value = '''8.39'
$ws.Range('E43').Value = '  -0.99%  '

$ws.Range('D44').Value = '2.858.77'
$ws.Range('E44').Value = '  -0.27%  '

$ws.Range('D45').Value = '''2.32'
$ws.Range('E45').Value = '  -3.66%  '

$ws.Range('D46').Value = '''0.269'
$ws.Range('E46').Value = '  -1.84%  '

$ws.Range('D47').Value = '''2.47'
$ws.Range('E47').Value = '  +5.81%  '

$ws.Range('D48').Value = '''26.59'
$ws.Range('E48').Value = '  -1.62%  '

$ws.Range('D49').Value = '''0.999'
$ws.Range('E49').Value = '  +0.03%  '

$ws.Range('E50').Value = '  -1.41%  '

$ws.Range('D51').Value = '''119.95'
$ws.Range('E51').Value = '  +0.11%  '
